$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.370.56'
$ws.Range("E2").Value = '  +0.66%  '

# Row 3
$ws.Range("D3").Value = '3.345.68'
$ws.Range("E3").Value = '  +1.51%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '189.40'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +4.26%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '558.76'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.20%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Value = '3.339.14'
$ws.Range("E8").Value = '  +1.48%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.584'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.06%  '

# Row 10
$ws.Range("E10").Value = '  -2.69%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.584'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.68%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '46.77'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.06%  '

# Row 13
$ws.Range("E13").Value = '  +1.69%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '8.61'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.44%  '

# Row 15
$ws.Range("D15").Value = '3.880.58'
$ws.Range("E15").Value = '  +1.46%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '598.53'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -6.34%  '

# Row 17
$ws.Range("D17").Value = '66.468.00'
$ws.Range("E17").Value = '  +0.82%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '18.01'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.47%  '

# Row 19
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.349.58'
$ws.Range("E19").Value = '  +1.53%  '

# Row 20
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.118'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.99%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '11.05'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.43%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.902'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.57%  '

# Row 23
$ws.Range("E23").Value = '  +4.21%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.04'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.26%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '99.67'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -7.39%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.97'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.50%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '5.98'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.54%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.38%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.58'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.47%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '8.69'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.36%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '30.59'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.63%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.73'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +5.83%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.94'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.78%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '581.33'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +4.20%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '11.02'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.70%  '

# Row 36
$ws.Range("E36").Value = '  -0.49%  '

# Row 37
$ws.Range("D37").Value = '3.710.60'
$ws.Range("E37").Value = '  +0.62%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.02%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '56.44'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.59%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.52'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.82%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '33.80'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +4.96%  '

# Row 42
$ws.Range("D42").Value = '0.0₃0707'
$ws.Range("E42").Value = '  -1.10%  '

# Row 43
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.127'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.22%  '

# Row 44
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.20'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -8.44%  '

# Row 45
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.66'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.52%  '

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.342'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.34%  '

# Row 47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.38'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +4.49%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0420'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.81%  '

# Row 49
$ws.Range("E49").Value = '  -0.64%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.57'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.85%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.24%  '
